$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header block ----
$ws.Range("C4").Value = "ENSemenC"
$ws.Range("E4").Value = "date: 23/05/2025"
$ws.Range("B5").Value = "MEMBRES DE L'EQUIPE: TIEU Mattéo, GERARD Paul"

# ---- Task rows 7-21 : B (row #, literal now instead of formula), C (task), D (coder), E (share) ----

# Row 7
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "Modélisation des plantes et cycle de vie"
$ws.Range("D7").Value = "Binôme"
$ws.Range("E7").Value = "50/50"

# Row 8
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Représentation des terrains et espace de jeu"
$ws.Range("D8").Value = "Binôme"
$ws.Range("E8").Value = "50/50"

# Row 9
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "Système de score et de progression"
$ws.Range("D9").Value = "Binôme"
$ws.Range("E9").Value = "50/50"

# Row 10
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = "Modélisation interactions complexes"
$ws.Range("D10").Value = "Mattéo"
$ws.Range("E10").NumberFormat = "0%"
$ws.Range("E10").Value = 1

# Row 11
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = "Conception du flux principal du programme"
$ws.Range("D11").Value = "Mattéo"
$ws.Range("E11").NumberFormat = "0%"
$ws.Range("E11").Value = 1

# Row 12
$ws.Range("B12").Value = 6
$ws.Range("C12").Value = "Définition des conditions de fin de jeu"
$ws.Range("D12").Value = "Mattéo"
$ws.Range("E12").NumberFormat = "0%"
$ws.Range("E12").Value = 1

# Row 13
$ws.Range("B13").Value = 7
$ws.Range("C13").Value = "Validation des actions (logique métier)"
$ws.Range("D13").Value = "Mattéo"
$ws.Range("E13").NumberFormat = "0%"
$ws.Range("E13").Value = 1

# Row 14
$ws.Range("B14").Value = 8
$ws.Range("C14").Value = "Gestion de l'affichage et interface utilisateur"
$ws.Range("D14").Value = "Mattéo"
$ws.Range("E14").NumberFormat = "0%"
$ws.Range("E14").Value = 1

# Row 15
$ws.Range("B15").Value = 9
$ws.Range("C15").Value = "Gestion des messages d'information et d'alerte"
$ws.Range("D15").Value = "Paul"
$ws.Range("E15").NumberFormat = "0%"
$ws.Range("E15").Value = 1

# Row 16
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = "Affichage du score et des statistiques"
$ws.Range("D16").Value = "Paul"
$ws.Range("E16").NumberFormat = "0%"
$ws.Range("E16").Value = 1

# Row 17
$ws.Range("B17").Value = 11
$ws.Range("C17").Value = "Traitement des choix du joueur & validation"
$ws.Range("D17").Value = "Paul"
$ws.Range("E17").NumberFormat = "0%"
$ws.Range("E17").Value = 1

# Row 18
$ws.Range("B18").Value = 12
$ws.Range("C18").Value = "Gestion du temps et des événements "
$ws.Range("D18").Value = "Paul"
$ws.Range("E18").NumberFormat = "0%"
$ws.Range("E18").Value = 1

# Row 19
$ws.Range("B19").Value = 13
$ws.Range("C19").Value = "Gestion des modes de jeu"
$ws.Range("D19").Value = "Paul"
$ws.Range("E19").NumberFormat = "0%"
$ws.Range("E19").Value = 1

# Row 20
$ws.Range("B20").Value = 14
$ws.Range("C20").Value = "Système de simulation temporelle"
$ws.Range("D20").Value = "Paul"
$ws.Range("E20").NumberFormat = "0%"
$ws.Range("E20").Value = 1

# Row 21
$ws.Range("B21").Value = 15
$ws.Range("C21").Value = "Tests, Débogage et Intégration du code"
$ws.Range("D21").Value = "Binôme"
$ws.Range("E21").Value = "50/50"

# ---- Remarks block (rows 23-25) ----
$ws.Range("B23").Style = "Normal"
$ws.Range("B23").Value = "La modélisation initiale des éléments centraux (plantes, terrains) a été menée en étroite collaboration par le binôme."

$ws.Range("B24").Style = "Normal"
$ws.Range("B24").Value = "Le codage majoritairement en binôme a conduit à une forte imbrication des contributions sur certaines fonctionnalités."

$ws.Range("B25").Style = "Normal"
$ws.Range("B25").Value = "La documentation et la rédaction du rapport ont été des efforts collaboratifs constants."

# ---- View state: selection on C19 ----
$ws.Range("C19").Select()
